$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 888.8
$ws.Range("J17").Value = 898.7692
$ws.Range("L17").Value = 2696.3076
$ws.Range("N17").Value = -3032.3076

$ws.Range("H19").Value = 576.5
$ws.Range("I19").Value = 634.2857
$ws.Range("J19").Value = 518.7143
$ws.Range("K19").Value = 634.2857
$ws.Range("L19").Value = 518.7143
$ws.Range("M19").Value = -459.2857
$ws.Range("N19").Value = -868.7143

$ws.Range("H33").Value = 336
$ws.Range("I33").Value = 393
$ws.Range("J33").Value = 303.42856
$ws.Range("K33").Value = 393
$ws.Range("L33").Value = 303.42856
$ws.Range("M33").Value = -164
$ws.Range("N33").Value = -761.4285600000001

$ws.Range("H116").Value = 1944.8889
$ws.Range("I116").Value = 1834
$ws.Range("J116").Value = 2166.6667
$ws.Range("K116").Value = 1834
$ws.Range("L116").Value = 2166.6667
$ws.Range("M116").Value = 1608
$ws.Range("N116").Value = -9050.6667

$ws.Range("H129").Value = 895.4026
$ws.Range("I129").Value = 766.0833
$ws.Range("K129").Value = 2298.2499
$ws.Range("M129").Value = 2701.7501

$ws.Range("H138").Value = 1501.7797
$ws.Range("I138").Value = 852.0454999999999
$ws.Range("J138").Value = 3407.6667
$ws.Range("K138").Value = 2556.1365
$ws.Range("L138").Value = 10223.0001
$ws.Range("M138").Value = 2583.8635
$ws.Range("N138").Value = -20503.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1463.1428
$ws.Range("I102").Value = 846
$ws.Range("J102").Value = 3006
$ws.Range("K102").Value = 846
$ws.Range("L102").Value = 3006
$ws.Range("M102").Value = 776
$ws.Range("N102").Value = -6250

$ws.Range("H132").Value = 17947064
$ws.Range("I132").Value = 19608492
$ws.Range("J132").Value = 7355466
$ws.Range("K132").Value = 58825476
$ws.Range("L132").Value = 22066398
$ws.Range("M132").Value = -58822946
$ws.Range("N132").Value = -22071458

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 32666.666
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 47500
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 47500
$ws.Range("M30").Value = -2875
$ws.Range("N30").Value = -47750

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H105").Value = 58825356
$ws.Range("I105").Value = 1778.3846
$ws.Range("J105").Value = 250001970
$ws.Range("K105").Value = 1778.3846
$ws.Range("L105").Value = 250001970
$ws.Range("M105").Value = -31.38460000000009
$ws.Range("N105").Value = -250005464

$ws.Range("H134").Value = 3178429.5
$ws.Range("I134").Value = 735.8946999999999
$ws.Range("J134").Value = 6951941
$ws.Range("K134").Value = 2207.6841
$ws.Range("L134").Value = 20855823
$ws.Range("M134").Value = 327.3159000000001
$ws.Range("N134").Value = -20860893

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 39800
$ws.Range("J64").Value = 39800
$ws.Range("L64").Value = 39800
$ws.Range("N64").Value = -40296

$ws.Range("H67").Value = 39800
$ws.Range("J67").Value = 39800
$ws.Range("L67").Value = 39800
$ws.Range("N67").Value = -41516

$ws.Range("H132").Value = 10102180
$ws.Range("I132").Value = 1249.9048
$ws.Range("J132").Value = 27778806
$ws.Range("K132").Value = 3749.7144
$ws.Range("L132").Value = 83336418
$ws.Range("M132").Value = -1219.7144
$ws.Range("N132").Value = -83341478

$ws.Range("H134").Value = 1485.9706
$ws.Range("I134").Value = 1338.1578
$ws.Range("J134").Value = 1673.2
$ws.Range("K134").Value = 4014.4734
$ws.Range("L134").Value = 5019.6
$ws.Range("M134").Value = -1479.4734
$ws.Range("N134").Value = -10089.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 27780694
$ws.Range("J5").Value = 27784126
$ws.Range("L5").Value = 83352378
$ws.Range("N5").Value = -83352602

$ws.Range("H86").Value = 466.5
$ws.Range("I86").Value = 146.4
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 439.2
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = 746.8
$ws.Range("N86").Value = -5372

$ws.Range("H89").Value = 466.5
$ws.Range("I89").Value = 146.4
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 1317.6
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = 4610.4
$ws.Range("N89").Value = -20856

$ws.Range("H115").Value = 8546.625
$ws.Range("I115").Value = 438.5
$ws.Range("J115").Value = 11249.333
$ws.Range("K115").Value = 1315.5
$ws.Range("L115").Value = 33747.999
$ws.Range("M115").Value = -140.5
$ws.Range("N115").Value = -36097.999

$ws.Range("H122").Value = 9474233
$ws.Range("I122").Value = 62500240
$ws.Range("K122").Value = 562502160
$ws.Range("M122").Value = -562499710

$ws.Range("H131").Value = 913.28
$ws.Range("I131").Value = 544.2222
$ws.Range("J131").Value = 949.7802
$ws.Range("K131").Value = 1632.6666
$ws.Range("L131").Value = 2849.3406
$ws.Range("M131").Value = 3407.3334
$ws.Range("N131").Value = -12929.3406

$ws.Range("H135").Value = 27780694
$ws.Range("J135").Value = 27784126
$ws.Range("L135").Value = 250057134
$ws.Range("N135").Value = -250062204

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 51.363636
$ws.Range("I2").Value = 26.666666
$ws.Range("J2").Value = 81
$ws.Range("K2").Value = 26.666666
$ws.Range("L2").Value = 81
$ws.Range("M2").Value = 86.33333400000001
$ws.Range("N2").Value = -307

$ws.Range("H132").Value = 4567.38
$ws.Range("I132").Value = 2895.7
$ws.Range("J132").Value = 11254.1
$ws.Range("K132").Value = 8687.099999999999
$ws.Range("L132").Value = 33762.3
$ws.Range("M132").Value = -6157.099999999999
$ws.Range("N132").Value = -38822.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1002.7857
$ws.Range("I81").Value = 894.0909
$ws.Range("J81").Value = 1401.3334
$ws.Range("K81").Value = 1788.1818
$ws.Range("L81").Value = 2802.6668
$ws.Range("M81").Value = -727.1818000000001
$ws.Range("N81").Value = -4924.6668

$ws.Range("H84").Value = 1002.7857
$ws.Range("I84").Value = 894.0909
$ws.Range("J84").Value = 1401.3334
$ws.Range("K84").Value = 8940.909
$ws.Range("L84").Value = 14013.334
$ws.Range("M84").Value = -3636.909
$ws.Range("N84").Value = -24621.334

$ws.Range("H122").Value = 41055.383
$ws.Range("I122").Value = 85216.664
$ws.Range("J122").Value = 3202.8572
$ws.Range("K122").Value = 255649.992
$ws.Range("L122").Value = 9608.571599999999
$ws.Range("M122").Value = -253199.992
$ws.Range("N122").Value = -14508.5716

$ws.Range("H132").Value = 8495539
$ws.Range("I132").Value = 31917.883
$ws.Range("K132").Value = 95753.649
$ws.Range("M132").Value = -93223.649

$ws.Range("H136").Value = 3362.3013
$ws.Range("I136").Value = 4472.4565
$ws.Range("J136").Value = 1470.9259
$ws.Range("K136").Value = 13417.3695
$ws.Range("L136").Value = 4412.7777
$ws.Range("M136").Value = -10867.3695
$ws.Range("N136").Value = -9512.777699999999
